$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

$ws.Range("G12").Value = "1/2/3/4:3STEP TX"
$ws.Range("G16").Value = "01:清償" + [char]10 + "02:法拍" + [char]10 + "03:增貸" + [char]10 + "04:展期" + [char]10 + "05:撥款" + [char]10 + "06:查閱" + [char]10 + "07:重估" + [char]10 + "08:其他"
$ws.Range("G22").Value = "1:重要袋" + [char]10 + "2:普通袋" + [char]10 + "3:重要袋&普通袋" + [char]10 + "4:其他"

$ws.Range("G23").Select()
